$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row before row 29, shifting existing rows 29+ down by one.
$ws.Rows.Item(29).Insert()

# Populate the new row 29 with the new September entry.
$ws.Cells.Item(29, 18).Value = "login internet invalid"
$ws.Cells.Item(29, 19).Value = "2024-09-03 19:54:49"
